# FA230_TestData_CreateOrMaintainAssetCategories_21C.xlsx
# Remove the hard-coded Oracle Cloud login URL / username / password that were
# stored (and hyperlinked) in the Input_Value sheet, cells AE2:AG2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input_Value")
$ws.Activate()

# Drop the hyperlink that lived on AE2 (pointed at the Oracle Cloud URL).
$ws.Range("AE2").Hyperlinks.Delete()

# Clear the credential values themselves (URL, user name, password).
$ws.Range("AE2:AG2").ClearContents()

# Put the cursor/selection back on the now-empty credential cells, scrolled
# over so they are visible (mirrors the saved view of the edited workbook).
[void]$ws.Range("AE2:AG2").Select()
$excel.ActiveWindow.ScrollColumn = 26
$excel.ActiveWindow.ScrollRow = 1
$null
